{"js": "// Add hierarchical section numbers (e.g. \"1 \", \"2.1 \", \"2.4.2 \") in front\n// of the document's Heading 2/3/4 titles. Numbering starts at the first\n// heading AFTER the \"Abstract\" section (i.e. \"Abstract\" itself stays\n// unnumbered, matching the site's chapter numbering scheme), and resets\n// sub-levels whenever a higher level advances.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text,items/style\");\nawait context.sync();\n\nconst LEVEL_STYLES = {\n  \"Heading 2\": 0,\n  \"Heading 3\": 1,\n  \"Heading 4\": 2,\n};\n\nlet numberingStarted = false;\nconst counters = [0, 0, 0];\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const paragraph = paragraphs.items[i];\n  const style = paragraph.style;\n\n  if (!numberingStarted) {\n    // \"Abstract\" is the last unnumbered heading; numbering begins right\n    // after it.\n    if (style === \"Heading 2\" && paragraph.text === \"Abstract\") {\n      numberingStarted = true;\n    }\n    continue;\n  }\n\n  const level = LEVEL_STYLES[style];\n  if (level === undefined) {\n    continue;\n  }\n\n  counters[level] += 1;\n  for (let deeper = level + 1; deeper < counters.length; deeper++) {\n    counters[deeper] = 0;\n  }\n  const sectionNumber = counters.slice(0, level + 1).join(\".\");\n\n  paragraph.insertText(sectionNumber + \" \", Word.InsertLocation.start);\n}\n\nawait context.sync();\n", "ps1": "# Add hierarchical section numbers (e.g. \"1 \", \"2.1 \", \"2.4.2 \") in front\n# of the document's Heading 2/3/4 titles. Numbering starts at the first\n# heading AFTER the \"Abstract\" section (i.e. \"Abstract\" itself stays\n# unnumbered, matching the site's chapter numbering scheme), and resets\n# sub-levels whenever a higher level advances.\n\n$d = $word.ActiveDocument\n\n$numberingStarted = $false\n$counters = @(0, 0, 0)\n\nforeach ($p in $d.Paragraphs) {\n    $styleName = $p.Style.NameLocal\n    $text = $p.Range.Text.Trim()\n\n    if (-not $numberingStarted) {\n        if ($styleName -eq \"Heading 2\" -and $text -eq \"Abstract\") {\n            $numberingStarted = $true\n        }\n        continue\n    }\n\n    $level = -1\n    if ($styleName -eq \"Heading 2\") { $level = 0 }\n    elseif ($styleName -eq \"Heading 3\") { $level = 1 }\n    elseif ($styleName -eq \"Heading 4\") { $level = 2 }\n\n    if ($level -eq -1) {\n        continue\n    }\n\n    $counters[$level] = $counters[$level] + 1\n    for ($deeper = $level + 1; $deeper -lt $counters.Length; $deeper++) {\n        $counters[$deeper] = 0\n    }\n\n    $parts = @()\n    for ($i = 0; $i -le $level; $i++) {\n        $parts += $counters[$i]\n    }\n    $sectionNumber = [string]::Join(\".\", $parts)\n\n    $p.Range.InsertBefore($sectionNumber + \" \")\n}\n"}
